$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Conexión de la base de datos con la aplicación web"
#    -> append " a través de C#"
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Conexión de la base de datos con la aplicación web",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Conexión de la base de datos con la aplicación web a través de C#", 2) | Out-Null

# ------------------------------------------------------------------
# 2) "Scripts de creación de base de datos" -> "Base de datos con Posgrest"
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Scripts de creación de base de datos",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Base de datos con Posgrest", 2) | Out-Null

# ------------------------------------------------------------------
# 3) "Vista de status de las oficinas" -> append " (Empleados)"
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Vista de status de las oficinas",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Vista de status de las oficinas (Empleados)", 2) | Out-Null

# ------------------------------------------------------------------
# 4) "Visualización de rutas (Empleados)" -> "Vista de rutas (Empleados)"
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Visualización de rutas (Empleados)",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Vista de rutas (Empleados)", 2) | Out-Null

# ------------------------------------------------------------------
# 5) Drop the two "Google Maps" bullets and replace them with two new
#    bullets: "Vista de rutas" and "Modificación de rutas (Empleados)"
#    (the latter carries the relocated _GoBack bookmark).
# ------------------------------------------------------------------

# Locate the two paragraphs by their current text.
$mapsPara1 = $null
$mapsPara2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -eq "Vista de rutas por medio de Google Maps`r") {
        $mapsPara1 = $i
    }
    if ($t -eq "Vista de ubicación de las oficinas por medio de Google Maps (Empleados)`r") {
        $mapsPara2 = $i
    }
}

$startRange = $d.Paragraphs.Item($mapsPara1).Range.Start
$endRange = $d.Paragraphs.Item($mapsPara2).Range.Start + $d.Paragraphs.Item($mapsPara2).Range.End - $d.Paragraphs.Item($mapsPara2).Range.Start
$endRange = $d.Paragraphs.Item($mapsPara2 + 1).Range.Start

$killRange = $d.Range($startRange, $endRange)
$killRange.Delete()

# The paragraph right before $startRange is now "Vista de tablas con DataTables".
$prevPara = $d.Paragraphs.Item($mapsPara1 - 1)
$prevPara.Range.InsertParagraphAfter()

$newPara1 = $d.Paragraphs.Item($mapsPara1)
$newPara1.Range.Text = "Vista de rutas"

$newPara1.Range.InsertParagraphAfter()

$newPara2 = $d.Paragraphs.Item($mapsPara1 + 1)
$newPara2.Range.Text = "Modificación de rutas (Empleados)"

# Re-anchor the _GoBack bookmark between "Modificación de " and "rutas (Empleados)".
$bmPoint = $newPara2.Range.Start + ("Modificación de ").Length
$bmRange = $d.Range($bmPoint, $bmPoint)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
